$wb = $excel.ActiveWorkbook

# "Pop" is the active sheet in the source workbook. Adding a worksheet
# inserts it immediately before whichever sheet is active and copies that
# sheet's formatting/content, so switch to "PIB" first to avoid cloning
# "Pop", then relocate the new sheet to the very end (after "Pop") so it
# becomes physically sheet5.xml and the last tab, matching the target.
$wb.Worksheets.Item("PIB").Activate()

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "PIB_change"
$newSheet.Move($null, $wb.Worksheets.Item("Pop"))

# Re-fetch the sheet reference by name: the COM object returned by Add()
# stops reflecting writes once the sheet has been relocated with Move().
$newSheet = $wb.Worksheets.Item("PIB_change")

$colA = @("mnemonico", "CodIBGE", "PIB_2016", "PIB_2017", "PIBCap2016", "PIBCap2017", "ChgPIB", "ChgPIBCap", "PopEst")
$colB = @("descricao", "codigo.ibge", "PIB a preços correntes de 2016", "PIB a preços correntes de 2017", "PIB per capta a preços correntes de 2016", "PIB per capta a preços correntes de 2017", "Variação do PIB de 2017 para 2016", "Variação do PIB per Capta de 2017 para 2016", "Populção Estimada")

for ($i = 0; $i -lt $colA.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 2).Value = $colB[$i]
}

# PIB sheet selection moves from A9 to B3.
$wb.Worksheets.Item("PIB").Range("B3").Select()

# Pop sheet: selection changes to A1:B1.
$wb.Worksheets.Item("Pop").Range("A1:B1").Select()

# New sheet becomes the active tab/selected sheet, with its own selection.
$newSheet.Range("B10").Select()
$newSheet.Activate()
